$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.183.42'
$ws.Range('E2').Value = '  -1.49%  '

$ws.Range('D3').Value = '2.996.86'
$ws.Range('E3').Value = '  -1.83%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '''586.43'
$ws.Range('E5').Value = '  -0.07%  '

$ws.Range('D6').Value = '''145.53'
$ws.Range('E6').Value = '  -3.77%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '''0.525'
$ws.Range('E8').Value = '  -2.19%  '

$ws.Range('D9').Value = '2.994.14'
$ws.Range('E9').Value = '  -1.95%  '

$ws.Range('E10').Value = '  -4.21%  '

$ws.Range('D11').Value = '''5.77'
$ws.Range('E11').Value = '  -0.81%  '

$ws.Range('D12').Value = '''0.463'
$ws.Range('E12').Value = '  +3.42%  '

$ws.Range('E13').Value = '  -2.49%  '

$ws.Range('D14').Value = '''34.44'
$ws.Range('E14').Value = '  -4.94%  '

$ws.Range('E15').Value = '  +1.88%  '

$ws.Range('D16').Value = '3.495.53'
$ws.Range('E16').Value = '  -1.68%  '

$ws.Range('D17').Value = '''7.05'
$ws.Range('E17').Value = '  -1.32%  '

$ws.Range('D18').Value = '62.113.89'
$ws.Range('E18').Value = '  -1.53%  '

$ws.Range('D19').Value = '2.994.10'
$ws.Range('E19').Value = '  -1.85%  '

$ws.Range('D20').Value = '''456.18'
$ws.Range('E20').Value = '  -4.30%  '

$ws.Range('D21').Value = '''13.96'
$ws.Range('E21').Value = '  -2.21%  '

$ws.Range('D22').Value = '''0.687'
$ws.Range('E22').Value = '  -2.53%  '

$ws.Range('D23').Value = '''7.39'
$ws.Range('E23').Value = '  -1.64%  '

$ws.Range('D24').Value = '''81.70'
$ws.Range('E24').Value = '  -0.75%  '

$ws.Range('E25').Value = '  -8.99%  '

$ws.Range('D26').Value = '''12.20'
$ws.Range('E26').Value = '  -3.80%  '

$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('D28').Value = '''9.73'
$ws.Range('E28').Value = '  -8.34%  '

$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.16%  '

$ws.Range('E30').Value = '  -1.55%  '

$ws.Range('D31').Value = '''6.94'
$ws.Range('E31').Value = '  -5.44%  '

$ws.Range('D32').Value = '''2.09'
$ws.Range('E32').Value = '  -4.87%  '

$ws.Range('D33').Value = '''27.62'
$ws.Range('E33').Value = '  -0.10%  '

$ws.Range('E34').Value = '  -1.75%  '

$ws.Range('D35').Value = '0.0₃0803'
$ws.Range('E35').Value = '  -1.79%  '

$ws.Range('D36').Value = '''1.02'
$ws.Range('E36').Value = '  -3.34%  '

$ws.Range('D37').Value = '''5.73'
$ws.Range('E37').Value = '  -2.67%  '

$ws.Range('D38').Value = '''2.10'
$ws.Range('E38').Value = '  -5.37%  '

$ws.Range('B39').Value = 'Cosmos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D39').Value = '''9.17'
$ws.Range('E39').Value = '  -0.58%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '''50.23'
$ws.Range('E40').Value = '  -0.47%  '

$ws.Range('E41').Value = '  +7.10%  '

$ws.Range('D42').Value = '''2.87'
$ws.Range('E42').Value = '  -11.44%  '

$ws.Range('D43').Value = '''391.46'
$ws.Range('E43').Value = '  -9.63%  '

$ws.Range('D44').Value = '''0.0357'
$ws.Range('E44').Value = '  -1.40%  '

$ws.Range('D45').Value = '''0.267'
$ws.Range('E45').Value = '  -7.19%  '

$ws.Range('D46').Value = '2.728.31'
$ws.Range('E46').Value = '  -3.53%  '

$ws.Range('D47').Value = '''37.36'
$ws.Range('E47').Value = '  -2.69%  '

$ws.Range('D48').Value = '''129.45'
$ws.Range('E48').Value = '  +0.13%  '

$ws.Range('D50').Value = '''0.109'
$ws.Range('E50').Value = '  -0.67%  '

$ws.Range('D51').Value = '''2.18'
$ws.Range('E51').Value = '  -0.83%  '
